$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume refresh - values are stored as text (inline strings)
# in the original sheet, so force Text number format before assignment to avoid
# Excel auto-converting numeric-looking strings (e.g. "211.27") into floating point numbers.
$updates = @{
    'D2' = '27.829.99'
    'E2' = '  -0.56%  '
    'D3' = '1.628.77'
    'E3' = '  -0.26%  '
    'E4' = '  -0.11%  '
    'D5' = '211.27'
    'E7' = '  -0.05%  '
    'D8' = '23.22'
    'E8' = '  -1.25%  '
    'E9' = '  -0.59%  '
    'E10' = '  -1.03%  '
    'E11' = '  -0.04%  '
    'D12' = '1.860.10'
    'E12' = '  -0.16%  '
    'D13' = '1.621.46'
    'E13' = '  -0.52%  '
    'D14' = '4.02'
    'E14' = '  -1.12%  '
    'E15' = '  -1.31%  '
    'D16' = '64.91'
    'E16' = '  -1.14%  '
    'D17' = '27.841.10'
    'E17' = '  -0.47%  '
    'D18' = '228.35'
    'E18' = '  -1.65%  '
    'E19' = '  +1.04%  '
    'E20' = '  -1.26%  '
    'E21' = '  -0.04%  '
    'D22' = '4.35'
    'E22' = '  -0.35%  '
    'E23' = '  -5.47%  '
    'E24' = '  -0.21%  '
    'D25' = '155.46'
    'E25' = '  +0.73%  '
    'E26' = '  -0.23%  '
    'D28' = '15.46'
    'E28' = '  -1.30%  '
    'D29' = '0.998'
    'E29' = '  +0.01%  '
    'E30' = '  -0.41%  '
    'E31' = '  -0.23%  '
    'E32' = '  -0.24%  '
    'E33' = '  -0.01%  '
    'D34' = '1.412.44'
    'E34' = '  +0.25%  '
    'E35' = '  +2.58%  '
    'E36' = '  -2.69%  '
    'E37' = '  -1.35%  '
    'E38' = '  -1.06%  '
    'E39' = '  -0.62%  '
    'E40' = '  -2.07%  '
    'E41' = '  -0.09%  '
    'E42' = '  -1.98%  '
    'D43' = '65.73'
    'E43' = '  -1.76%  '
    'B44' = 'FraxShare'
    'C44' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D44' = '5.43'
    'E44' = '  -0.81%  '
    'B45' = 'RenderToken'
    'C45' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D45' = '1.81'
    'E45' = '  -0.28%  '
    'D46' = '1.769.19'
    'E46' = '  -0.26%  '
    'E48' = '  +0.33%  '
    'E49' = '  +0.96%  '
    'D50' = '0.0503'
    'E50' = '  -0.37%  '
    'D51' = '7.62'
    'E51' = '  +0.66%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
